$d = $word.ActiveDocument

function Find-ParagraphStartingWith($prefix) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

function Set-ParagraphRunXml($paragraph, $runsXml) {
    $full = $paragraph.Range
    # Exclude the trailing paragraph mark so pPr / numbering stay intact.
    $r = $d.Range($full.Start, $full.End - 1)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $runsXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# --- "Dispositivo (IDDisp, Tipo, NomePianta^REPLICA , NumReplica^REPLICA) "
#   -> "Dispositivo (IDDisp, IDReplica^REPLICA ,Tipo)"
# Now "(IDDisp, IDReplica^REPLICA " is underlined instead of just "IDDisp";
# the red colouring on NomePianta/NumReplica is removed together with those runs.
$dispositivoRuns = '<w:r><w:rPr><w:b/></w:rPr><w:t>Dispositivo</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>(IDDisp,</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>IDReplica</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>REPLICA</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t>,</w:t></w:r>' +
  '<w:r><w:t>Tipo</w:t></w:r>' +
  '<w:r><w:t>)</w:t></w:r>'

$pDispositivo = Find-ParagraphStartingWith "Dispositivo"
if ($pDispositivo -ne $null) {
    Set-ParagraphRunXml $pDispositivo $dispositivoRuns
}

# --- "Replica (NumReplica, NomePianta^SPECIE, Gruppo, DataDimora, Esposizione)"
#   -> "Replica (IDReplica, Gruppo, DataDimora, Esposizione, NomePianta^SPECIE)"
# "NumReplica" (red+underlined) becomes "IDReplica" (underlined only, no colour);
# " NomePianta" keeps the red colour on the leading space only, and the
# superscript "SPECIE" loses its colour/underline formatting.
$replicaRuns = '<w:r><w:rPr><w:b/></w:rPr><w:t>Replica</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>ID</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Replica</w:t></w:r>' +
  '<w:r><w:t>,</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> Gruppo, DataDimora, Esposizione</w:t></w:r>' +
  '<w:r><w:t>,</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t>NomePianta</w:t></w:r>' +
  '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>SPECIE</w:t></w:r>' +
  '<w:r><w:t>)</w:t></w:r>'

$pReplica = Find-ParagraphStartingWith "Replica "
if ($pReplica -ne $null) {
    Set-ParagraphRunXml $pReplica $replicaRuns
}

Write-Output "Done."
